# Auto-assisted generated script - Generate Report for Handoff
# Adds two new tracked files (438c604b..., ea1c7e31...) into the
# localization-status report across Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

# ===================== Sheet1: Overview =====================
$ws1 = $wb.Worksheets.Item(1)

# Clear existing hyperlinks on the sheet so we can rebuild them cleanly
$ws1.Cells.Hyperlinks.Delete()

# Insert a row above the current row 3 (for new file 438c604b) and a new
# row at the end (row 5, for new file ea1c7e31)
$ws1.Range("3:3").Insert()
$ws1.Range("5:5").Insert()

# Row 3: 438c604b-....md
$ws1.Range("A3").Value = "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-30-19 10:30:08"

# Row 5: ea1c7e31-....md
$ws1.Range("A5").Value = "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-30-19 10:30:08"

# Rebuild hyperlinks for column A (rows 2-5)
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e567af052fa7de912a8cc3f2c2f4f5e35b597cb1/e2e/dc0e6672-dfa2-4dd3-b629-b2e8897fedcb.md", "", "", "dc0e6672-dfa2-4dd3-b629-b2e8897fedcb.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/e2e/438c604b-9f4f-419a-ba3f-e6ef65b44f61.md", "", "", "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0256384309014c2d47bd3141ba47df580f521432/e2e/6ca15e1d-8973-4ab8-888d-185a31d3e070.md", "", "", "6ca15e1d-8973-4ab8-888d-185a31d3e070.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/e2e/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md", "", "", "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md") | Out-Null

# ===================== ws2: zh-cn =====================
$ws2 = $wb.Worksheets.Item(2)

# Clear existing hyperlinks on the sheet so we can rebuild them cleanly
$ws2.Cells.Hyperlinks.Delete()

# Insert a row above current row 3 (for 438c604b) and a new row at the
# end (row 5, for ea1c7e31)
$ws2.Range("3:3").Insert()
$ws2.Range("5:5").Insert()

# Row 3: 438c604b-9f4f-419a-ba3f-e6ef65b44f61
$ws2.Range("A3").Value = "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-19 10:30:05"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

# Row 4: 6ca15e1d-8973-4ab8-888d-185a31d3e070
$ws2.Range("A4").Value = "6ca15e1d-8973-4ab8-888d-185a31d3e070.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-19 10:28:51"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "Include"

# Row 5: ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2
$ws2.Range("A5").Value = "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-19 10:30:05"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("I5").Value = "Include"

# Rebuild hyperlinks
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/e2e/438c604b-9f4f-419a-ba3f-e6ef65b44f61.md", "", "", "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/e2e/438c604b-9f4f-419a-ba3f-e6ef65b44f61.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.zh-cn.xlf", "", "", "438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0256384309014c2d47bd3141ba47df580f521432/e2e/6ca15e1d-8973-4ab8-888d-185a31d3e070.md", "", "", "6ca15e1d-8973-4ab8-888d-185a31d3e070.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/0256384309014c2d47bd3141ba47df580f521432/e2e/6ca15e1d-8973-4ab8-888d-185a31d3e070.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c31a5832dcad8e63a456a635ecea11642893eaa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.zh-cn.xlf", "", "", "6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/e2e/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md", "", "", "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/e2e/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md", "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.zh-cn.xlf", "", "", "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.zh-cn.xlf") | Out-Null

# ===================== ws3: de-de =====================
$ws3 = $wb.Worksheets.Item(3)

# Clear existing hyperlinks on the sheet so we can rebuild them cleanly
$ws3.Cells.Hyperlinks.Delete()

# Insert a row above current row 3 (for 438c604b) and a new row at the
# end (row 5, for ea1c7e31)
$ws3.Range("3:3").Insert()
$ws3.Range("5:5").Insert()

# Row 3: 438c604b-9f4f-419a-ba3f-e6ef65b44f61
$ws3.Range("A3").Value = "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-19 10:30:08"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

# Row 4: 6ca15e1d-8973-4ab8-888d-185a31d3e070
$ws3.Range("A4").Value = "6ca15e1d-8973-4ab8-888d-185a31d3e070.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-19 10:28:54"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "Include"

# Row 5: ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2
$ws3.Range("A5").Value = "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-19 10:30:08"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("I5").Value = "Include"

# Rebuild hyperlinks
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/e2e/438c604b-9f4f-419a-ba3f-e6ef65b44f61.md", "", "", "438c604b-9f4f-419a-ba3f-e6ef65b44f61.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/e2e/438c604b-9f4f-419a-ba3f-e6ef65b44f61.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/438c604b9f4f419aba3fe6ef65b44f6100000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.de-de.xlf", "", "", "438c604b-9f4f-419a-ba3f-e6ef65b44f61.d691c3bf86f9f2204045abe490a600dd0c0ea3bf.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0256384309014c2d47bd3141ba47df580f521432/e2e/6ca15e1d-8973-4ab8-888d-185a31d3e070.md", "", "", "6ca15e1d-8973-4ab8-888d-185a31d3e070.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/0256384309014c2d47bd3141ba47df580f521432/e2e/6ca15e1d-8973-4ab8-888d-185a31d3e070.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/613ad9ecaa392dd515cfbfed4c1f74e4c210f0de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.de-de.xlf", "", "", "6ca15e1d-8973-4ab8-888d-185a31d3e070.d89eab08ad5bd30c317853cd37cba76e1827ebb5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/e2e/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md", "", "", "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/e2e/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.md", "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea1c7e312e604e37acd38f3c8f4262e200000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.de-de.xlf", "", "", "ea1c7e31-2e60-4e37-acd3-8f3c8f4262e2.467d0b36da4103842d8b15e314d3118fc76ca854.de-de.xlf") | Out-Null

Write-Host "Done"
